$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K10").Value = 4000
$ws.Range("M10").Value = 6631.67
$ws.Range("N10").Value = 5164.03
$ws.Range("O10").Value = 4236.55

$ws.Range("K11").Value = 14482.07
$ws.Range("M11").Value = 827429.11
$ws.Range("N11").Value = 479988.61
$ws.Range("O11").Value = 428252.01

$ws.Range("M12").Value = 137737.98
$ws.Range("N12").Value = 58380.77
$ws.Range("O12").Value = 55897.04

$ws.Range("K13").Value = 1300
$ws.Range("M13").Value = 21272.16
$ws.Range("N13").Value = 15074.72
$ws.Range("O13").Value = 13323.32

$ws.Range("N14").Value = 5416.22
$ws.Range("O14").Value = 5416.22

$ws.Range("K17").Value = 39136.81

$ws.Range("O19").Value = 3256.05

$ws.Range("K21").Value = 1394.92
$ws.Range("M21").Value = 2297.25

$ws.Range("K22").Value = 1396.67
$ws.Range("M22").Value = 1933.45

$ws.Range("N26").Value = 70210
$ws.Range("O26").Value = 69770
